$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0278
$ws.Range("G2").Value = -0.08226950354609927
$ws.Range("H2").Value = -0.08226950354609927
$ws.Range("I2").Value = -0.2671394799054373
$ws.Range("J2").Value = -0.2671394799054373
$ws.Range("K2").Value = -1.62
$ws.Range("L2").Value = -0.3829787234042553
$ws.Range("U2").Value = 1.21
$ws.Range("V2").Value = 0.01704225352112676
$ws.Range("W2").Value = -0.1664953751284687
$ws.Range("X2").Value = 0.08209565768076103
$ws.Range("Y2").Value = -0.2485910328092297
$ws.Range("Z2").Value = 0.1003463491009157
$ws.Range("AA2").Value = -0.02680647150922807
$ws.Range("AB2").Value = 0.06334338498171473
$ws.Range("AC2").Value = -0.09014985649094279
$ws.Range("AD2").Value = 37.7
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 37.7
$ws.Range("AG2").Value = 36.49
$ws.Range("AH2").Value = 0.3468261269549218
$ws.Range("AI2").Value = 0.8231441048034934
$ws.Range("AJ2").Value = 0.3394734393897106
$ws.Range("AK2").Value = 0.8183449203857367
$ws.Range("AL2").Value = 0.718
$ws.Range("AM2").Value = 0.718
$ws.Range("AN2").Value = 100
$ws.Range("AO2").Value = -1.573816155988858
$ws.Range("AP2").Value = 96.79045092838197
$ws.Range("AQ2").Value = -1.573816155988858

# Row 3
$ws.Range("D3").Value = 0.0278
$ws.Range("G3").Value = -0.08226950354609927
$ws.Range("H3").Value = -0.08226950354609927
$ws.Range("I3").Value = -0.2671394799054373
$ws.Range("J3").Value = -0.2671394799054373
$ws.Range("K3").Value = -1.62
$ws.Range("L3").Value = -0.3829787234042553
$ws.Range("U3").Value = 1.21
$ws.Range("V3").Value = 0.01704225352112676
$ws.Range("W3").Value = -0.1664953751284687
$ws.Range("X3").Value = 0.08209565768076103
$ws.Range("Y3").Value = -0.2485910328092297
$ws.Range("Z3").Value = 0.1003463491009157
$ws.Range("AA3").Value = -0.02680647150922807
$ws.Range("AB3").Value = 0.06334338498171473
$ws.Range("AC3").Value = -0.09014985649094279
$ws.Range("AD3").Value = 37.7
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 37.7
$ws.Range("AG3").Value = 36.49
$ws.Range("AH3").Value = 0.3468261269549218
$ws.Range("AI3").Value = 0.8231441048034934
$ws.Range("AJ3").Value = 0.3394734393897106
$ws.Range("AK3").Value = 0.8183449203857367
$ws.Range("AL3").Value = 0.718
$ws.Range("AM3").Value = 0.718
$ws.Range("AN3").Value = 100
$ws.Range("AO3").Value = -1.573816155988858
$ws.Range("AP3").Value = 96.79045092838197
$ws.Range("AQ3").Value = -1.573816155988858

